# Apply the edits described by the commit:
#  - Relabel several row headers in column A (rows 34-36, 41-44)
#  - Add new survey data in columns S:W for rows 25-44
#  - Add the per-row average formula in column X for rows 25-44
#  - Update the sheet view (zoom + selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update row labels (column A) that were reworded ---
$ws.Range("A34").Value = "Durchschnittliche Verweildauer p. Bild (4)"
$ws.Range("A35").Value = "Durchschnittliche Verweildauer p. B. (20)"
$ws.Range("A36").Value = "Durchschnittliche Verweildauer p.B. (100)"
$ws.Range("A41").Value = "Fixation Dauer pro Bild 1"
$ws.Range("A42").Value = "Fixation Dauer pro Bild 4"
$ws.Range("A43").Value = "Fixation Dauer pro Bild 20"
$ws.Range("A44").Value = "Fixation Dauer pro Bild 100"

# --- 2. Fill new data columns S:W and add the running-average formula in X ---
#     (columns S-W hold 5 more participants' measurements; X = SUM(D:W)/20)
# Row 25
$ws.Range("S25").Value = 6
$ws.Range("T25").Value = 30
$ws.Range("U25").Value = 5
$ws.Range("V25").Value = 34
$ws.Range("W25").Value = 11
$ws.Range("X25").Formula = "=SUM(D25:W25)/20"

# Row 26
$ws.Range("S26").Value = 4
$ws.Range("T26").Value = 4
$ws.Range("U26").Value = 99
$ws.Range("V26").Value = 95
$ws.Range("W26").Value = 4
$ws.Range("X26").Formula = "=SUM(D26:W26)/20"

# Row 27
$ws.Range("S27").Value = 18
$ws.Range("T27").Value = 32
$ws.Range("U27").Value = 52
$ws.Range("V27").Value = 16
$ws.Range("W27").Value = 20
$ws.Range("X27").Formula = "=SUM(D27:W27)/20"

# Row 28
$ws.Range("S28").Value = 65
$ws.Range("T28").Value = 53
$ws.Range("U28").Value = 41
$ws.Range("V28").Value = 77
$ws.Range("W28").Value = 74
$ws.Range("X28").Formula = "=SUM(D28:W28)/20"

# Row 29
$ws.Range("S29").Value = 3.8333333333333335
$ws.Range("T29").Value = 1.5333333333333334
$ws.Range("U29").Value = 2.8
$ws.Range("V29").Value = 2.8823529411764706
$ws.Range("W29").Value = 3.1818181818181817
$ws.Range("X29").Formula = "=SUM(D29:W29)/20"

# Row 30
$ws.Range("S30").Value = 10.75
$ws.Range("T30").Value = 6
$ws.Range("U30").Value = 2.9393939393939394
$ws.Range("V30").Value = 2.221052631578947
$ws.Range("W30").Value = 4.5
$ws.Range("X30").Formula = "=SUM(D30:W30)/20"

# Row 31
$ws.Range("S31").Value = 7.944444444444445
$ws.Range("T31").Value = 2.8125
$ws.Range("U31").Value = 4.365384615384615
$ws.Range("V31").Value = 7.6875
$ws.Range("W31").Value = 4.15
$ws.Range("X31").Formula = "=SUM(D31:W31)/20"

# Row 32
$ws.Range("S32").Value = 9.261538461538462
$ws.Range("T32").Value = 2.5849056603773586
$ws.Range("U32").Value = 3
$ws.Range("V32").Value = 4.090909090909091
$ws.Range("W32").Value = 6.72972972972973
$ws.Range("X32").Formula = "=SUM(D32:W32)/20"

# Row 33
$ws.Range("S33").Value = 3369.5
$ws.Range("T33").Value = 1697.1666666666667
$ws.Range("U33").Value = 2010.6
$ws.Range("V33").Value = 2206.9411764705883
$ws.Range("W33").Value = 3257.4545454545455
$ws.Range("X33").Formula = "=SUM(D33:W33)/20"

# Row 34
$ws.Range("S34").Value = 3167.75
$ws.Range("T34").Value = 1839.5
$ws.Range("U34").Value = 1129
$ws.Range("V34").Value = 850.5368421052632
$ws.Range("W34").Value = 2108.75
$ws.Range("X34").Formula = "=SUM(D34:W34)/20"

# Row 35
$ws.Range("S35").Value = 1205.4444444444443
$ws.Range("T35").Value = 633.0625
$ws.Range("U35").Value = 876.8269230769231
$ws.Range("V35").Value = 2399.4375
$ws.Range("W35").Value = 979.7
$ws.Range("X35").Formula = "=SUM(D35:W35)/20"

# Row 36
$ws.Range("S36").Value = 1058.553846153846
$ws.Range("T36").Value = 381.41509433962267
$ws.Range("U36").Value = 316.3170731707317
$ws.Range("V36").Value = 661.5454545454545
$ws.Range("W36").Value = 912.8243243243244
$ws.Range("X36").Formula = "=SUM(D36:W36)/20"

# Row 37
$ws.Range("S37").Value = 4.5
$ws.Range("T37").Value = 2.1
$ws.Range("U37").Value = 2.8
$ws.Range("V37").Value = 5
$ws.Range("W37").Value = 4.2727272727272725
$ws.Range("X37").Formula = "=SUM(D37:W37)/20"

# Row 38
$ws.Range("S38").Value = 4.5
$ws.Range("T38").Value = 3.5
$ws.Range("U38").Value = 2.4646464646464645
$ws.Range("V38").Value = 1.8631578947368421
$ws.Range("W38").Value = 3.75
$ws.Range("X38").Formula = "=SUM(D38:W38)/20"

# Row 39
$ws.Range("S39").Value = 2.611111111111111
$ws.Range("T39").Value = 1.34375
$ws.Range("U39").Value = 1.5576923076923077
$ws.Range("V39").Value = 4.4375
$ws.Range("W39").Value = 1.85
$ws.Range("X39").Formula = "=SUM(D39:W39)/20"

# Row 40
$ws.Range("S40").Value = 2.0153846153846153
$ws.Range("T40").Value = 0.8113207547169812
$ws.Range("U40").Value = 0.34146341463414637
$ws.Range("V40").Value = 1.2987012987012987
$ws.Range("W40").Value = 1.5675675675675675
$ws.Range("X40").Formula = "=SUM(D40:W40)/20"

# Row 41
$ws.Range("S41").Value = 2527.1666666666665
$ws.Range("T41").Value = 1564.5
$ws.Range("U41").Value = 1566.8
$ws.Range("V41").Value = 1841
$ws.Range("W41").Value = 2683.4545454545455
$ws.Range("X41").Formula = "=SUM(D41:W41)/20"

# Row 42
$ws.Range("S42").Value = 2081.75
$ws.Range("T42").Value = 1174.25
$ws.Range("U42").Value = 869.969696969697
$ws.Range("V42").Value = 682.7052631578947
$ws.Range("W42").Value = 1587.5
$ws.Range("X42").Formula = "=SUM(D42:W42)/20"

# Row 43
$ws.Range("S43").Value = 801.6666666666666
$ws.Range("T43").Value = 503.75
$ws.Range("U43").Value = 548.0576923076923
$ws.Range("V43").Value = 1947.875
$ws.Range("W43").Value = 687.45
$ws.Range("X43").Formula = "=SUM(D43:W43)/20"

# Row 44
$ws.Range("S44").Value = 680.8153846153846
$ws.Range("T44").Value = 270.52830188679246
$ws.Range("U44").Value = 118.1951219512195
$ws.Range("V44").Value = 488.2597402597403
$ws.Range("W44").Value = 628.081081081081
$ws.Range("X44").Formula = "=SUM(D44:W44)/20"

# --- 3. Update the sheet view: new zoom level and active selection ---
$ws.Range("Y43").Select()
$window = $excel.ActiveWindow
$window.Zoom = 83

